# Expand the "hasLinkTo" test data in properties.xlsx:
#   - add a new row (21) for "linkstoRegion" (Region)
#   - fill in the previously-blank rows 22-24 and complete rows 25-26 with
#     hasLinkTo* rows for: StillImageRepresentation, Resource,
#     ArchiveRepresentation, MovingImageRepresentation, AudioRepresentation
#   - update the sheet view (zoom level + current selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet view: zoom + selection ---
$excel.ActiveWindow.Zoom = 75
$ws.Range("N26").Select()

# --- Row 21 (new row): linkstoRegion / Region -------------------------
$ws.Range("A21").Value = "linkstoRegion"
$ws.Range("B21").Value = "hasLinkTo"
$ws.Range("C21").Value = "Region"
$ws.Range("D21").Value = "Links to a region of an image"
$ws.Range("E21").Value = "Verweist auf eine Region in einem Bild"
$ws.Range("F21").Value = "Se réfère à une région d’une image"
$ws.Range("N21").Value = "Searchbox"
$ws.Range("O21").Value = "numprops: 1"

# --- Row 22: hasLinkToImage / StillImageRepresentation -----------------
$ws.Range("A22").Value = "hasLinkToImage"
$ws.Range("B22").Value = "hasLinkTo"
$ws.Range("C22").Value = "StillImageRepresentation"
$ws.Range("D22").Value = "link to image"
$ws.Range("N22").Value = "Searchbox"
# newly-touched blank cells keep the row's existing style (1), not the
# column default (5)
$ws.Range("A22").Copy()
$ws.Range("I22:M22").PasteSpecial(-4122)

# --- Row 23: hasLinkToResource / Resource -------------------------------
$ws.Range("A23").Value = "hasLinkToResource"
$ws.Range("B23").Value = "hasLinkTo"
$ws.Range("C23").Value = "Resource"
$ws.Range("N23").Value = "Searchbox"
$ws.Range("D23").Value = "hasLinkToResource"
$ws.Range("A23").Copy()
$ws.Range("D23").PasteSpecial(-4122)

# --- Row 24: hasLinkToArchiveRepresentation / ArchiveRepresentation -----
$ws.Range("A24").Value = "hasLinkToArchiveRepresentation"
$ws.Range("B24").Value = "hasLinkTo"
$ws.Range("C24").Value = "ArchiveRepresentation"
$ws.Range("N24").Value = "Searchbox"
$ws.Range("D24").Value = "hasLinkToArchiveRepresentation"
$ws.Range("A24").Copy()
$ws.Range("D24").PasteSpecial(-4122)

# --- Row 25: hasLinkToMovingImageRepesentation / MovingImageRepresentation
$ws.Range("A25").Value = "hasLinkToMovingImageRepesentation"
$ws.Range("B25").Value = "hasLinkTo"
$ws.Range("C25").Value = "MovingImageRepresentation"
$ws.Range("D25").Value = "hasLinkToMovingImageRepesentation"
$ws.Range("N25").Value = "Searchbox"
$ws.Range("A25").Copy()
$ws.Range("B25:D25").PasteSpecial(-4122)
$ws.Range("N25").PasteSpecial(-4122)

# --- Row 26: hasLinkToAudioRepesentation / AudioRepresentation ----------
$ws.Range("A26").Value = "hasLinkToAudioRepesentation"
$ws.Range("B26").Value = "hasLinkTo"
$ws.Range("C26").Value = "AudioRepresentation"
$ws.Range("D26").Value = "hasLinkToAudioRepesentation"
$ws.Range("N26").Value = "Searchbox"
$ws.Range("A26").Copy()
$ws.Range("C26:D26").PasteSpecial(-4122)
